# Update COVID-19 Valais daily figures for rows 384-387 (2021-03-15 to 2021-03-18)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

function Set-NumericValue($cell, $value) {
    # Some cells (columns L/M) carry a Text ("@") number format. Writing a
    # plain numeric .Value into such a cell would be stored as a text
    # string, so we briefly switch to a numeric format, assign the value,
    # then restore the original (Text) format/style.
    $origFormat = $cell.NumberFormat
    $cell.NumberFormat = "General"
    $cell.Value = $value
    $cell.NumberFormat = $origFormat
}

# Row 384 (2021-03-15): new positive cases revised from 81 to 82
$ws.Range("C384").Value = 82

# Row 385 (2021-03-16): new positive cases revised from 60 to 80;
# one new in-hospital COVID-19 death recorded
$ws.Range("C385").Value = 80
Set-NumericValue $ws.Range("L385") 1

# Row 386 (2021-03-17): new positive cases revised from 19 to 82;
# one new in-hospital COVID-19 death recorded
$ws.Range("C386").Value = 82
Set-NumericValue $ws.Range("L386") 1

# Row 387 (2021-03-18): fill in the day's figures (previously blank)
$ws.Range("C387").Value = 11
$ws.Range("E387").Value = 5
$ws.Range("F387").Value = 4
$ws.Range("G387").Value = 27
Set-NumericValue $ws.Range("L387") 0
Set-NumericValue $ws.Range("M387") 0

$wb.Save()
